$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 11.467718355161836
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = 16.976148379153372
$ws.Range("E2").Value = $null

$ws.Range("B3").Value = 10.83486683656362
$ws.Range("C3").Value = -3.105531684919832
$ws.Range("D3").Value = 18.13267575692705
$ws.Range("E3").Value = -0.39689215022412583

$ws.Range("B1:E3").Select()
